$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 2.219956787131591
$ws.Range("C2").Value = 9.105099971597772
$ws.Range("D2").Value = -8.419520236579334
$ws.Range("E2").Value = -0.6983358980328844
$ws.Range("F2").Value = 0.8904915159445659
$ws.Range("G2").Value = -2.07292559449643
$ws.Range("H2").Value = 0.1925130480743557
$ws.Range("I2").Value = -0.2579532017493174
$ws.Range("J2").Value = -0.3979806818825675
$ws.Range("K2").Value = 0.1284229565523166

# Row 3
$ws.Range("B3").Value = 8.189548673647696
$ws.Range("C3").Value = -8.974959151229303
$ws.Range("D3").Value = -1.01725560823065
$ws.Range("E3").Value = 0.6812891314769711
$ws.Range("F3").Value = -2.226143026305237
$ws.Range("G3").Value = 0.06705817846833073
$ws.Range("H3").Value = -0.3695255935427733
$ws.Range("I3").Value = -0.5026284358164976
$ws.Range("J3").Value = 0.02723167798804321
$ws.Range("K3").Value = -0.1127918163358273

# Row 4
$ws.Range("B4").Value = -16.14018903381625
$ws.Range("C4").Value = -6.791910405777589
$ws.Range("D4").Value = -3.960768829455048
$ws.Range("E4").Value = -5.976670854937382
$ws.Range("F4").Value = -2.989026990346052
$ws.Range("G4").Value = -2.889230364707919
$ws.Range("H4").Value = -2.610106546763337
$ws.Range("I4").Value = -1.76447406392936
$ws.Range("J4").Value = -1.663120324238024
$ws.Range("K4").Value = -0.8349890726314544

# Row 5
$ws.Range("B5").Value = -5.111729259208365
$ws.Range("C5").Value = 3.455275933053657
$ws.Range("D5").Value = -3.728614538712288
$ws.Range("E5").Value = 0.9923340248431218
$ws.Range("F5").Value = -0.9110484627979438
$ws.Range("G5").Value = -0.1981003006930162
$ws.Range("H5").Value = -0.1744636854283967
$ws.Range("I5").Value = -0.02054777792822543
$ws.Range("J5").Value = 0.4483335670758678
$ws.Range("K5").Value = 0.6752115876730926

# Row 6
$ws.Range("B6").Value = -0.2719305936643379
$ws.Range("C6").Value = -1.865743883948665
$ws.Range("D6").Value = 0.2314251006721421
$ws.Range("E6").Value = -0.3608174671756919
$ws.Range("F6").Value = -0.4453535416162557
$ws.Range("G6").Value = 0.08612045519262451
$ws.Range("H6").Value = -0.06432879815791362
$ws.Range("I6").Value = 0.580021064545069
$ws.Range("J6").Value = 0.7052415095167848
$ws.Range("K6").Value = 0.3072362793633781

# Row 7
$ws.Range("B7").Value = -1.457511430326895
$ws.Range("C7").Value = 0.4830481445585361
$ws.Range("D7").Value = -0.4427235209988271
$ws.Range("E7").Value = -0.3702614922349701
$ws.Range("F7").Value = 0.1920276157063483
$ws.Range("G7").Value = -0.005912151999258208
$ws.Range("H7").Value = 0.6464597483751605
$ws.Range("I7").Value = 0.7799234681276213
$ws.Range("J7").Value = 0.3775495712643568
$ws.Range("K7").Value = 0.3009015042266514

# Row 8
$ws.Range("B8").Value = 0.377548932709721
$ws.Range("C8").Value = -0.568828804964185
$ws.Range("D8").Value = -0.2020134236017385
$ws.Range("E8").Value = 0.2526917083102306
$ws.Range("F8").Value = 0.005220687375374677
$ws.Range("G8").Value = 0.709322499702091
$ws.Range("H8").Value = 0.8372749959213646
$ws.Range("I8").Value = 0.4220547482687358
$ws.Range("J8").Value = 0.35197537525307
$ws.Range("K8").Value = 0.6305200013052115

# Row 9
$ws.Range("B9").Value = -1.173269424811924
$ws.Range("C9").Value = -0.2662168009564319
$ws.Range("D9").Value = 0.5442390163674474
$ws.Range("E9").Value = -0.03248618871055892
$ws.Range("F9").Value = 0.70073847962863
$ws.Range("G9").Value = 0.9210160735187328
$ws.Range("H9").Value = 0.4568818931967809
$ws.Range("I9").Value = 0.378539758305367
$ws.Range("J9").Value = 0.6760188420518104
$ws.Range("K9").Value = 0.4968561201928344

# Row 10
$ws.Range("B10").Value = -0.3538576297953431
$ws.Range("C10").Value = 0.4918198143377789
$ws.Range("D10").Value = 0.02616020827444671
$ws.Range("E10").Value = 0.7056676249657313
$ws.Range("F10").Value = 0.9099227904698376
$ws.Range("G10").Value = 0.4696485924218674
$ws.Range("H10").Value = 0.3872693173527261
$ws.Range("I10").Value = 0.6785207862006621
$ws.Range("J10").Value = 0.5031320576645155
$ws.Range("K10").Value = 0.6535938957934404

# Row 11
$ws.Range("B11").Value = 0.5142477939707875
$ws.Range("C11").Value = 0.02876374783689539
$ws.Range("D11").Value = 0.683599939311716
$ws.Range("E11").Value = 0.9039408101654827
$ws.Range("F11").Value = 0.4643053422394562
$ws.Range("G11").Value = 0.3769007839823598
$ws.Range("H11").Value = 0.6700434640770384
$ws.Range("I11").Value = 0.4953597493577047
$ws.Range("J11").Value = 0.6449699222604717
$ws.Range("K11").Value = 0.3699768206640927

# Row 12
$ws.Range("B12").Value = 0.05018831929411549
$ws.Range("C12").Value = 0.8037777881109285
$ws.Range("D12").Value = 0.8147464657516995
$ws.Range("E12").Value = 0.4345783383908784
$ws.Range("F12").Value = 0.3855269040122699
$ws.Range("G12").Value = 0.6441032906516032
$ws.Range("H12").Value = 0.4725811863294772
$ws.Range("I12").Value = 0.631351887111854
$ws.Range("J12").Value = 0.3514560234789169
$ws.Range("K12").Value = 0.1878533956106533

# Row 13
$ws.Range("B13").Value = 0.761457155790539
$ws.Range("C13").Value = 0.7855419918959663
$ws.Range("D13").Value = 0.4315577385859573
$ws.Range("E13").Value = 0.3666656204007895
$ws.Range("F13").Value = 0.623772159586399
$ws.Range("G13").Value = 0.4580510308952294
$ws.Range("H13").Value = 0.6147108416195859
$ws.Range("I13").Value = 0.3338721780501667
$ws.Range("J13").Value = 0.1713795120806105
$ws.Range("K13").Value = 0.04814754147425859

# Row 14
$ws.Range("B14").Value = 1.125739732513876
$ws.Range("C14").Value = 0.5069288003698943
$ws.Range("D14").Value = 0.1755723340877904
$ws.Range("E14").Value = 0.6484660721731321
$ws.Range("F14").Value = 0.4674304344515254
$ws.Range("G14").Value = 0.5599000430081948
$ws.Range("H14").Value = 0.3149780871614969
$ws.Range("I14").Value = 0.1565155112655648
$ws.Range("J14").Value = 0.01943583084488787
$ws.Range("K14").Value = 0.5403464745801891

# Row 15
$ws.Range("B15").Value = 0.9576236723601449
$ws.Range("C15").Value = 0.2226830363001488
$ws.Range("D15").Value = 0.4081828413823239
$ws.Range("E15").Value = 0.5033755155173882
$ws.Range("F15").Value = 0.557849285568602
$ws.Range("G15").Value = 0.2384098477254667
$ws.Range("H15").Value = 0.1291161168902814
$ws.Range("I15").Value = -0.006853219352837503
$ws.Range("J15").Value = 0.4969491838668565
$ws.Range("K15").Value = 0.2970525035592049

# Row 16
$ws.Range("B16").Value = 0.5350086232236873
$ws.Range("C16").Value = 0.5414715811178482
$ws.Range("D16").Value = 0.3242970937011966
$ws.Range("E16").Value = 0.5860435243391693
$ws.Range("F16").Value = 0.2744351589178602
$ws.Range("G16").Value = 0.09157454668606776
$ws.Range("H16").Value = -0.01239555844019874
$ws.Range("I16").Value = 0.5011245128056051
$ws.Range("J16").Value = 0.2858677898194339

# Row 17
$ws.Range("B17").Value = 0.7770620722256893
$ws.Range("C17").Value = 0.4089758345755062
$ws.Range("D17").Value = 0.4404431577054521
$ws.Range("E17").Value = 0.2843889017771028
$ws.Range("F17").Value = 0.1069397163856127
$ws.Range("G17").Value = -0.05124389339762181
$ws.Range("H17").Value = 0.4852787037784192
$ws.Range("I17").Value = 0.2775335613519331

# Row 18
$ws.Range("B18").Value = 0.720000021056889
$ws.Range("C18").Value = 0.5575319321166473
$ws.Range("D18").Value = 0.121924760327536
$ws.Range("E18").Value = 0.1358986668508143
$ws.Range("F18").Value = -0.01550102094271821
$ws.Range("G18").Value = 0.4539510573947921
$ws.Range("H18").Value = 0.2743085116504074

# Row 19
$ws.Range("B19").Value = 0.8043571095207618
$ws.Range("C19").Value = 0.1395806278654255
$ws.Range("D19").Value = 0.04127642773981979
$ws.Range("E19").Value = 0.01671636730576187
$ws.Range("F19").Value = 0.4663391832225094
$ws.Range("G19").Value = 0.2534447081011285

# Row 20
$ws.Range("B20").Value = 0.3792569872816374
$ws.Range("C20").Value = 0.1263006712090333
$ws.Range("D20").Value = -0.09976821105225223
$ws.Range("E20").Value = 0.4814444548743619
$ws.Range("F20").Value = 0.2766837437271186

# Row 21
$ws.Range("B21").Value = 0.2911607482431205
$ws.Range("C21").Value = -0.08623093880433602
$ws.Range("D21").Value = 0.4184715358843989
$ws.Range("E21").Value = 0.2867219094086165

# Row 22
$ws.Range("B22").Value = 0.167618103061072
$ws.Range("C22").Value = 0.5177895860664353
$ws.Range("D22").Value = 0.1751453671933744

# Row 23
$ws.Range("B23").Value = 0.5618492773058843
$ws.Range("C23").Value = 0.1965658720679752

# Row 24
$ws.Range("B24").Value = 0.4328090033804217
